# Apply conscientiousness (and related) score corrections to surveydata sheet
# The survey scoring values in columns V, W, X, Y, Z for several respondents were
# recomputed (halved) for the "conscientiousness" column (W) on most rows, and for
# row 38 the V, W, Y, Z scores were also corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "W2" = 2
    "W3" = 1
    "W4" = 3
    "W5" = 1
    "W6" = 4
    "W7" = 2
    "W8" = 4
    "W9" = 3
    "W10" = 1
    "W11" = 3
    "W12" = 1
    "W13" = 4
    "W14" = 4
    "W15" = 3
    "W16" = 4
    "W17" = 1
    "W18" = 2
    "W19" = 2
    "W20" = 4
    "W21" = 4
    "W22" = 1
    "W23" = 2
    "W25" = 1
    "W26" = 3
    "W27" = 2
    "W28" = 3
    "W29" = 4
    "W30" = 4
    "W31" = 4
    "W32" = 3
    "W33" = 2
    "W34" = 3
    "W35" = 5
    "W36" = 3
    "W37" = 2
    "V38" = 2
    "W38" = 2
    "Y38" = 1
    "Z38" = 2
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

